# Fruta / hortaliza, semanal
# A new weekly observation is inserted at row 199 (Terminal Hortofrutícola
# Agro Chillán - Mango), pushing the former rows 199:217 down to 200:218.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 199; Excel shifts rows 199:217 down to
# 200:218 and the sheet dimension grows from A1:T217 to A1:T218.
$ws.Rows.Item(199).Insert()

# Fill in the new row with the latest weekly record.
$ws.Range("A199").Value = 7
$ws.Range("B199").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C199").Value = "Ñuble"
$ws.Range("D199").Value = 45265
$ws.Range("E199").Value = 16
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100108
$ws.Range("H199").Value = "Tropicales y subtropicales"
$ws.Range("I199").Value = 100108002
$ws.Range("J199").Value = "Mango"
$ws.Range("K199").Value = "Sin especificar"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 100
$ws.Range("N199").Value = 12000
$ws.Range("O199").Value = 12000
$ws.Range("P199").Value = 12000
$ws.Range("Q199").Value = "$/bandeja 4 kilos"
$ws.Range("R199").Value = "Perú"
$ws.Range("S199").Value = 3000
$ws.Range("T199").Value = 4
